$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AA2").Value = 27
$ws.Range("AE2").Value = 24
$ws.Range("AO2").Value = 18.5
$ws.Range("G2").Value = 4.1
$ws.Range("H2").Value = 2.1
$ws.Range("J2").Value = 3.4
$ws.Range("K2").Value = 3.6
$ws.Range("P2").Value = 1.82
$ws.Range("W2").Value = 1.32
$ws.Range("F3").Value = 1.48
$ws.Range("G3").Value = 1.61
$ws.Range("T3").Value = 1.87
$ws.Range("W3").Value = 2.6
$ws.Range("AB4").Value = 9.199999999999999
$ws.Range("O4").Value = 1.31
$ws.Range("S4").Value = 3.4
$ws.Range("AE5").Value = 290
$ws.Range("AM5").Value = 300
$ws.Range("H5").Value = 13
$ws.Range("I5").Value = 14
$ws.Range("J5").Value = 5.5
$ws.Range("V5").Value = 1.07
$ws.Range("F6").Value = 2.5
$ws.Range("I6").Value = 3
$ws.Range("Q6").Value = 1.82
$ws.Range("T6").Value = 1.59
$ws.Range("U6").Value = 1.94
$ws.Range("V6").Value = 1.5
$ws.Range("N7").Value = 4.2
$ws.Range("R7").Value = 1.44
$ws.Range("U7").Value = 2
$ws.Range("N8").Value = 5.3
$ws.Range("U8").Value = 2.12
$ws.Range("N10").Value = 5.4
$ws.Range("T10").Value = 1.65
$ws.Range("AN11").Value = 8.800000000000001
$ws.Range("N11").Value = 5.3
$ws.Range("X11").Value = 21
$ws.Range("AF12").Value = 130
$ws.Range("F12").Value = 12
$ws.Range("P12").Value = 4.3
$ws.Range("S12").Value = 1.7
$ws.Range("T12").Value = 1.65
$ws.Range("W12").Value = 1.08
$ws.Range("H13").Value = 29
$ws.Range("I13").Value = 32
$ws.Range("R13").Value = 1.68
$ws.Range("W13").Value = 7.4
$ws.Range("AB14").Value = 7.6
$ws.Range("AE14").Value = 200
$ws.Range("AL14").Value = 44
$ws.Range("H14").Value = 9.199999999999999
$ws.Range("J14").Value = 4.8
$ws.Range("N14").Value = 3.9
$ws.Range("Q14").Value = 1.84
$ws.Range("S14").Value = 3.25
$ws.Range("T14").Value = 2.12
$ws.Range("U14").Value = 1.78
$ws.Range("X14").Value = 16.5
$ws.Range("Y14").Value = 970
$ws.Range("AB15").Value = 980
$ws.Range("AH15").Value = 22
$ws.Range("AN15").Value = 95
$ws.Range("N15").Value = 5.2
$ws.Range("P15").Value = 2.44
$ws.Range("Q15").Value = 1.56
$ws.Range("R15").Value = 1.57
$ws.Range("U15").Value = 2.08
$ws.Range("F16").Value = 1.04
$ws.Range("G16").Value = 1000
$ws.Range("H16").Value = 1.04
$ws.Range("I16").Value = 1000
$ws.Range("J16").Value = 1.01
$ws.Range("K16").Value = 980
$ws.Range("M16").Value = 1.02
$ws.Range("P16").Value = 1.25
$ws.Range("R16").Value = 1.2
$ws.Range("S16").Value = 1.59
$ws.Range("T16").Value = 1.01
$ws.Range("U16").Value = 1.89
$ws.Range("V16").Value = 1.01
$ws.Range("W16").Value = 1.01
$ws.Range("F17").Value = 2.66
$ws.Range("G17").Value = 2.98
$ws.Range("I17").Value = 2.98
$ws.Range("J17").Value = 3.3
$ws.Range("K17").Value = 3.9
$ws.Range("M17").Value = 1.05
$ws.Range("N17").Value = 3.85
$ws.Range("O17").Value = 1.3
$ws.Range("P17").Value = 1.97
$ws.Range("R17").Value = 1.34
$ws.Range("S17").Value = 2.92
$ws.Range("V17").Value = 1.5
$ws.Range("W17").Value = 1.5
$ws.Range("AC18").Value = 17.5
$ws.Range("AD18").Value = 46
$ws.Range("AF18").Value = 1000
$ws.Range("AH18").Value = 34
$ws.Range("H18").Value = 5
$ws.Range("I18").Value = 13
$ws.Range("K18").Value = 7.6
$ws.Range("V18").Value = 1.08
$ws.Range("Y18").Value = 46
$ws.Range("Z18").Value = 100
